$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the existing last row's timestamp.
$ws.Range("A19").Value = 45865.83358738426

# Append the new row of sensor data recorded by the scheduled task.
$ws.Range("A20").Value = 45865.87529501117
$ws.Range("A20").NumberFormat = $ws.Range("A19").NumberFormat
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 13.83
$ws.Range("E20").Value = 90.41
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 8.34
$ws.Range("H20").Value = "SE"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = "21:00:25"
